# Regenerate save_data: use K (strikeouts) values instead of the old Strike# values.
# Column G holds the "K" stat per the header in row 1 (G1 = "K").
# The new K values below were recalculated from the source data and are written
# as literal numbers into column G for each data row.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$kValues = @{
    2  = 0
    3  = 2
    4  = 1
    5  = 0
    6  = 1
    7  = 4
    8  = 1
    9  = 0
    10 = 2
    11 = 1
    12 = 2
    13 = 1
    14 = 0
    16 = 2
    17 = 2
    18 = 0
    19 = 0
    20 = 0
    21 = 2
    22 = 1
    23 = 1
    24 = 0
    25 = 0
    26 = 1
    27 = 3
    29 = 1
    30 = 0
    32 = 0
    33 = 2
    34 = 0
    35 = 1
    36 = 0
    38 = 0
    40 = 1
    42 = 0
    43 = 1
    44 = 1
    45 = 1
}

foreach ($row in $kValues.Keys) {
    $ws.Cells.Item($row, 7).Value = $kValues[$row]
}
